# Update cryptos list prices/volumes (GitHub Actions scrape refresh).
# Column D ("Price") and column E ("Volume(1h)") are stored as plain text
# in the sheet, so numeric-looking prices must be forced to text
# (NumberFormat "@") before assignment to avoid Excel silently
# re-interpreting them as numbers and dropping trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.911.69"
$ws.Range("E2").Value = "  -3.88%  "
$ws.Range("D3").Value = "1.633.84"
$ws.Range("E3").Value = "  -6.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.54"
$ws.Range("E5").Value = "  -4.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4720"
$ws.Range("E7").Value = "  -6.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2547"
$ws.Range("E8").Value = "  -6.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06024"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07017"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").Value = "1.638.54"
$ws.Range("E11").Value = "  -6.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.72"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6134"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.351"
$ws.Range("E14").Value = "  -6.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.50"
$ws.Range("E15").Value = "  -6.67%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9982"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "24.923.03"
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006527"
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("E20").Value = "  -6.52%  "
$ws.Range("D21").Value = "1.844.95"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.365"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.561"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.252"
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.37"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.74"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.368"
$ws.Range("E27").Value = "  -8.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "102.62"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.640"
$ws.Range("E29").Value = "  -7.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.749"
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07716"
$ws.Range("E31").Value = "  -6.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.542"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9987"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04283"
$ws.Range("E34").Value = "  -8.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.594"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9185"
$ws.Range("E36").Value = "  -7.36%  "
$ws.Range("E37").Value = "  -6.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.558"
$ws.Range("E38").Value = "  -6.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01543"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9979"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8230"
$ws.Range("E41").Value = "  +8.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.788"
$ws.Range("E42").Value = "  -6.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.88"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3691"
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.715"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1097"
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05213"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.045"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.43"
$ws.Range("E49").Value = "  -3.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9992"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.22%  "
